# Refactored and renamed method DecrementPosition() to Demote()
#
# This updates the "Rules: Moving Down in the Queue" worked example to use
# the new Demote/Promote terminology, and tweaks a label in the "Remove"
# example to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Text updates in the "Demote (Moving Down in the Queue)" example (rows 11-18) ---
$ws.Range("J12").Value = "Rules: Demote (Moving Down in the Queue)"
$ws.Range("E14").Value = "Demote B to P.5"
$ws.Range("J14").Value = "Promote items that are lower than the current position and higher than or equal to the requested position"
$ws.Range("J15").Value = "No change for items lower in the queue than requested position"

# The old sub-heading cell ("Move Down") is no longer used in the refactored example
$ws.Range("E17").ClearContents()

# --- Text update in the "Removing a specific item" example (rows 21-28) ---
$ws.Range("J24").Value = "Promote items lower in the queue than the item being removed "

# --- Formatting / view tweaks ---
# Column E needs to widen slightly to fit the new "Demote"/"Promote" wording
# (ColumnWidth is in characters; 15.1 lands on the OOXML serialized width of 16).
$ws.Columns.Item(5).ColumnWidth = 15.1

# The worked example now centres on row 14 ("Demote B to P.5"), so move the
# selection there (matches the scrolled view saved in the workbook).
$ws.Range("J14").Select()
